$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Player")
$ws.Range("F1").Value = "TestLongDescription"
$ws.Range("F1").Font.Name = "Arial"
$ws.Range("F1").Font.Size = 10
$ws.Range("F1").Font.Color = 0
